$wb = $excel.ActiveWorkbook

function Add-Link {
    param($ws, $addr, $url, $disp)
    $ws.Hyperlinks.Add($ws.Range($addr), $url, [Type]::Missing, [Type]::Missing, $disp) | Out-Null
}

# ---------------------------------------------------------------------------
# Sheet "Overview": drop the two now-obsolete ".png" dependency rows and
# replace the remaining source-file row with the two newly generated
# markdown source files. The ".localization-config" row is unchanged in
# content but shifts up to row 4 once the obsolete row is removed.
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

# Remove the old row 4 (b9643a17-....png, IsDependency) so the
# ".localization-config" row shifts from row 5 up to row 4.
$wsOverview.Rows.Item(4).Delete()

# Rebuild hyperlinks cleanly (stale references would otherwise survive).
$wsOverview.Hyperlinks.Delete()

$wsOverview.Range("B2").Value = "Ready for handoff"
$wsOverview.Range("C2").Value = "Ready for handoff"
Add-Link $wsOverview "A2" "https://github.com/OpenLocalizationTest/oltest/blob/69155895447040a691a4c7c5f9c3eacc8275ed43/e2e/5d1200b2-a8b8-48f0-aba7-73ccccb6f9dc.md" "5d1200b2-a8b8-48f0-aba7-73ccccb6f9dc.md"

$wsOverview.Range("B3").Value = "Ready for handoff"
$wsOverview.Range("C3").Value = "Ready for handoff"
Add-Link $wsOverview "A3" "https://github.com/OpenLocalizationTest/oltest/blob/69155895447040a691a4c7c5f9c3eacc8275ed43/e2e/fdbf1f95-06fc-4b80-9c8d-929ce072bed6.md" "fdbf1f95-06fc-4b80-9c8d-929ce072bed6.md"

$wsOverview.Range("B4").Value = "Not to be localized"
$wsOverview.Range("C4").Value = "Not to be localized"
Add-Link $wsOverview "A4" "https://github.com/OpenLocalizationTest/oltest/blob/69155895447040a691a4c7c5f9c3eacc8275ed43/.localization-config" ".localization-config"

# ---------------------------------------------------------------------------
# "zh-cn" detail sheet
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

# Remove the old IsDependency row (row 4) so the ".localization-config"
# row shifts up from row 5 to row 4, keeping its content intact.
$ws.Rows.Item(4).Delete()
$ws.Hyperlinks.Delete()

$ws.Range("B2").Value = "Ready for handoff"
$ws.Range("D2").Value = "2016-03-09 01:31:48"
$ws.Range("E2").ClearContents()
$ws.Range("F2").ClearContents()
$ws.Range("G2").Value = "0001-01-01 00:00:00"
$ws.Range("H2").Value = "Include"
$ws.Range("I2").ClearContents()
Add-Link $ws "A2" "https://github.com/OpenLocalizationTest/oltest/blob/69155895447040a691a4c7c5f9c3eacc8275ed43/e2e/5d1200b2-a8b8-48f0-aba7-73ccccb6f9dc.md" "5d1200b2-a8b8-48f0-aba7-73ccccb6f9dc.md"
Add-Link $ws "C2" "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/c5d068c72a4ddfe6b02fa2b64bb2748d96b75e61/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/5d1200b2-a8b8-48f0-aba7-73ccccb6f9dc.b9b637c7fdad1bc415da70fe45ebd71a6a6a97d1.zh-cn.xlf" "5d1200b2-a8b8-48f0-aba7-73ccccb6f9dc.b9b637c7fdad1bc415da70fe45ebd71a6a6a97d1.zh-cn.xlf"

$ws.Range("B3").Value = "Ready for handoff"
$ws.Range("D3").Value = "2016-03-09 01:31:48"
$ws.Range("E3").ClearContents()
$ws.Range("F3").ClearContents()
$ws.Range("G3").Value = "0001-01-01 00:00:00"
$ws.Range("H3").Value = "Include"
$ws.Range("I3").ClearContents()
Add-Link $ws "A3" "https://github.com/OpenLocalizationTest/oltest/blob/69155895447040a691a4c7c5f9c3eacc8275ed43/e2e/fdbf1f95-06fc-4b80-9c8d-929ce072bed6.md" "fdbf1f95-06fc-4b80-9c8d-929ce072bed6.md"
Add-Link $ws "C3" "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/c5d068c72a4ddfe6b02fa2b64bb2748d96b75e61/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/fdbf1f95-06fc-4b80-9c8d-929ce072bed6.c3131bf1b5abe9ac1dc05897bab8b8956e9aa739.zh-cn.xlf" "fdbf1f95-06fc-4b80-9c8d-929ce072bed6.c3131bf1b5abe9ac1dc05897bab8b8956e9aa739.zh-cn.xlf"

# Row 4 (".localization-config") content already correct after the row
# shift; only the hyperlink needs to be (re)created.
Add-Link $ws "A4" "https://github.com/OpenLocalizationTest/oltest/blob/69155895447040a691a4c7c5f9c3eacc8275ed43/.localization-config" ".localization-config"

# ---------------------------------------------------------------------------
# "de-de" detail sheet
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("de-de")

$ws2.Rows.Item(4).Delete()
$ws2.Hyperlinks.Delete()

$ws2.Range("B2").Value = "Ready for handoff"
$ws2.Range("D2").Value = "2016-03-09 01:31:58"
$ws2.Range("E2").ClearContents()
$ws2.Range("F2").ClearContents()
$ws2.Range("G2").Value = "0001-01-01 00:00:00"
$ws2.Range("H2").Value = "Include"
$ws2.Range("I2").ClearContents()
Add-Link $ws2 "A2" "https://github.com/OpenLocalizationTest/oltest/blob/69155895447040a691a4c7c5f9c3eacc8275ed43/e2e/5d1200b2-a8b8-48f0-aba7-73ccccb6f9dc.md" "5d1200b2-a8b8-48f0-aba7-73ccccb6f9dc.md"
Add-Link $ws2 "C2" "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/3fa71b4c059f662ebea8f2068355744385aa2dbb/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/5d1200b2-a8b8-48f0-aba7-73ccccb6f9dc.b9b637c7fdad1bc415da70fe45ebd71a6a6a97d1.de-de.xlf" "5d1200b2-a8b8-48f0-aba7-73ccccb6f9dc.b9b637c7fdad1bc415da70fe45ebd71a6a6a97d1.de-de.xlf"

$ws2.Range("B3").Value = "Ready for handoff"
$ws2.Range("D3").Value = "2016-03-09 01:31:58"
$ws2.Range("E3").ClearContents()
$ws2.Range("F3").ClearContents()
$ws2.Range("G3").Value = "0001-01-01 00:00:00"
$ws2.Range("H3").Value = "Include"
$ws2.Range("I3").ClearContents()
Add-Link $ws2 "A3" "https://github.com/OpenLocalizationTest/oltest/blob/69155895447040a691a4c7c5f9c3eacc8275ed43/e2e/fdbf1f95-06fc-4b80-9c8d-929ce072bed6.md" "fdbf1f95-06fc-4b80-9c8d-929ce072bed6.md"
Add-Link $ws2 "C3" "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/3fa71b4c059f662ebea8f2068355744385aa2dbb/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/fdbf1f95-06fc-4b80-9c8d-929ce072bed6.c3131bf1b5abe9ac1dc05897bab8b8956e9aa739.de-de.xlf" "fdbf1f95-06fc-4b80-9c8d-929ce072bed6.c3131bf1b5abe9ac1dc05897bab8b8956e9aa739.de-de.xlf"

Add-Link $ws2 "A4" "https://github.com/OpenLocalizationTest/oltest/blob/69155895447040a691a4c7c5f9c3eacc8275ed43/.localization-config" ".localization-config"
